# Auto-generated Excel COM-interop script to apply the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link) --------------------------------------
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("B39").Value = 'Cronos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'

# --- Price column (D) -------------------------------------------------------
# Force text format per-cell (NumberFormat on a multi-area Range only sticks
# to the first area), so the dotted/locale-looking price strings aren't
# reinterpreted as numbers; then restore the default (unstyled) cell style.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.990.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.045.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.661'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.39'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.378'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0773'
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.62'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.344.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.805'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.054.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '36.983.29'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0893'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '235.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.65'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0609'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.42'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0895'
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.20'
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.106'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.15'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.22'
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '94.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.272.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.227.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.40'
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) column (E) ---------------------------------------------------
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +4.90%  '
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("E12").Value = '  +3.83%  '
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("E15").Value = '  +6.29%  '
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("E18").Value = '  +15.73%  '
$ws.Range("E19").Value = '  +3.48%  '
$ws.Range("E20").Value = '  -2.26%  '
$ws.Range("E21").Value = '  +0.87%  '
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  -2.25%  '
$ws.Range("E25").Value = '  +10.80%  '
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("E28").Value = '  -2.21%  '
$ws.Range("E29").Value = '  +0.79%  '
$ws.Range("E30").Value = '  +6.73%  '
$ws.Range("E31").Value = '  +2.31%  '
$ws.Range("E32").Value = '  -2.19%  '
$ws.Range("E33").Value = '  +2.53%  '
$ws.Range("E34").Value = '  +2.47%  '
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  +3.51%  '
$ws.Range("E40").Value = '  +13.31%  '
$ws.Range("E41").Value = '  +25.96%  '
$ws.Range("E42").Value = '  -1.07%  '
$ws.Range("E43").Value = '  -4.42%  '
$ws.Range("E44").Value = '  -1.23%  '
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("E46").Value = '  +2.78%  '
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("E48").Value = '  -1.92%  '
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("E50").Value = '  -1.57%  '
$ws.Range("E51").Value = '  -15.41%  '
